# Appended spectroscopy update: replace the single data row (row 2) of
# Sheet1 with the new subject's Subject/Scan_Date/Process_Date labels and
# refreshed ventilation-summary metrics, keeping row 1's headers and the
# overall A1:AB2 layout / formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: identifying text columns (A:C) -------------------------------
$ws.Range("A2").Value = "Xe-037"
$ws.Range("B2").Value = "2022-06-06"
$ws.Range("C2").Value = "2022-06-17"

# --- Row 2: numeric measurement columns (D:AB) ----------------------------
$ws.Range("D2").Value  = 18.40078544613359
$ws.Range("E2").Value  = 1.1210269486086077
$ws.Range("F2").Value  = 8.075790214094166
$ws.Range("G2").Value  = 0
$ws.Range("H2").Value  = 9.1968171627027733
$ws.Range("I2").Value  = 0.15166835187057634
$ws.Range("J2").Value  = 2.7256341495581835
$ws.Range("K2").Value  = 0
$ws.Range("L2").Value  = 2.8773025014287597
$ws.Range("M2").Value  = 3.9038114916252691
$ws.Range("N2").Value  = 9.6056622851365017
$ws.Range("O2").Value  = 23.680045720314766
$ws.Range("P2").Value  = 33.848419571811668
$ws.Range("Q2").Value  = 23.495405987602762
$ws.Range("R2").Value  = 5.4666549435090337
$ws.Range("S2").Value  = 2.1211588341319731
$ws.Range("T2").Value  = 3.5433243944256385
$ws.Range("U2").Value  = 9.0649316393370558
$ws.Range("V2").Value  = 35.800325317624306
$ws.Range("W2").Value  = 42.249527410207939
$ws.Range("X2").Value  = 7.2207324042730905
$ws.Range("Y2").Value  = 7.4929881337648325
$ws.Range("Z2").Value  = 14.727076591154262
$ws.Range("AA2").Value = 33.659115426105714
$ws.Range("AB2").Value = 44.120819848975188
